# Auto-generated: apply scheduled-runner market-data refresh to Sheets (per diff)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3187.8572
$ws.Range("I106").Value = 2863
$ws.Range("K106").Value = 2863
$ws.Range("M106").Value = -2232

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 8852.59
$ws.Range("I132").Value = 3007.0286
$ws.Range("K132").Value = 9021.085800000001
$ws.Range("M132").Value = -6491.085800000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2251.5588
$ws.Range("I137").Value = 2342.4814
$ws.Range("K137").Value = 7027.4442
$ws.Range("M137").Value = -4477.4442

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10315.833
$ws.Range("I32").Value = 10800
$ws.Range("J32").Value = 4990
$ws.Range("K32").Value = 10800
$ws.Range("L32").Value = 4990
$ws.Range("M32").Value = -10513
$ws.Range("N32").Value = -5564

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9109.857
$ws.Range("I61").Value = 4533.778
$ws.Range("J61").Value = 12541.917
$ws.Range("K61").Value = 4533.778
$ws.Range("L61").Value = 12541.917
$ws.Range("M61").Value = -4321.778
$ws.Range("N61").Value = -12965.917

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1904.3684
$ws.Range("I132").Value = 1500.3235
$ws.Range("J132").Value = 5338.75
$ws.Range("K132").Value = 4500.970499999999
$ws.Range("L132").Value = 16016.25
$ws.Range("M132").Value = -1970.970499999999
$ws.Range("N132").Value = -21076.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 9109.857
$ws.Range("I136").Value = 4533.778
$ws.Range("J136").Value = 12541.917
$ws.Range("K136").Value = 13601.334
$ws.Range("L136").Value = 37625.751
$ws.Range("M136").Value = -11051.334
$ws.Range("N136").Value = -42725.751

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 542.3333
$ws.Range("I22").Value = 505
$ws.Range("J22").Value = 617
$ws.Range("K22").Value = 505
$ws.Range("L22").Value = 617
$ws.Range("M22").Value = -155
$ws.Range("N22").Value = -1317

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2643.7083
$ws.Range("I31").Value = 1472.5
$ws.Range("K31").Value = 1472.5
$ws.Range("M31").Value = -1177.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2643.7083
$ws.Range("I34").Value = 1472.5
$ws.Range("K34").Value = 1472.5
$ws.Range("M34").Value = -1270.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2519.3
$ws.Range("I99").Value = 2410.3333
$ws.Range("J99").Value = 3500
$ws.Range("K99").Value = 2410.3333
$ws.Range("L99").Value = 3500
$ws.Range("M99").Value = -912.3332999999998
$ws.Range("N99").Value = -6496

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1472.45
$ws.Range("I107").Value = 450.33334
$ws.Range("K107").Value = 450.33334
$ws.Range("M107").Value = 1469.66666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2519.3
$ws.Range("I126").Value = 2410.3333
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 7230.999899999999
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -4760.999899999999
$ws.Range("N126").Value = -15440

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2826.0513
$ws.Range("I134").Value = 1341.3704
$ws.Range("J134").Value = 6166.5835
$ws.Range("K134").Value = 4024.1112
$ws.Range("L134").Value = 18499.7505
$ws.Range("M134").Value = -1489.1112
$ws.Range("N134").Value = -23569.7505

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 159226
$ws.Range("J141").Value = 159226
$ws.Range("L141").Value = 159226
$ws.Range("N141").Value = -169586

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6108.048
$ws.Range("I5").Value = 770.7778
$ws.Range("J5").Value = 10111
$ws.Range("K5").Value = 2312.3334
$ws.Range("L5").Value = 30333
$ws.Range("M5").Value = -2200.3334
$ws.Range("N5").Value = -30557

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 6810.6
$ws.Range("I22").Value = 184.33333
$ws.Range("J22").Value = 16750
$ws.Range("K22").Value = 552.99999
$ws.Range("L22").Value = 50250
$ws.Range("M22").Value = -383.99999
$ws.Range("N22").Value = -50588

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 6810.6
$ws.Range("I27").Value = 184.33333
$ws.Range("J27").Value = 16750
$ws.Range("K27").Value = 552.99999
$ws.Range("L27").Value = 50250
$ws.Range("M27").Value = -450.99999
$ws.Range("N27").Value = -50454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 9703.333000000001
$ws.Range("J130").Value = 11111
$ws.Range("L130").Value = 33333
$ws.Range("N130").Value = -43373

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1069.75
$ws.Range("J132").Value = 1100
$ws.Range("L132").Value = 9900
$ws.Range("N132").Value = -14960

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 6108.048
$ws.Range("I135").Value = 770.7778
$ws.Range("J135").Value = 10111
$ws.Range("K135").Value = 6937.000199999999
$ws.Range("L135").Value = 90999
$ws.Range("M135").Value = -4402.000199999999
$ws.Range("N135").Value = -96069

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2838.9
$ws.Range("I137").Value = 2292
$ws.Range("J137").Value = 3286.3635
$ws.Range("K137").Value = 6876
$ws.Range("L137").Value = 9859.0905
$ws.Range("M137").Value = -1776
$ws.Range("N137").Value = -20059.0905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 3550

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 30000
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 30000
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 30000
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -30504

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 186289.25
$ws.Range("J34").Value = 182000
$ws.Range("L34").Value = 182000
$ws.Range("N34").Value = -182536

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H76").Value = 186289.25
$ws.Range("J76").Value = 182000
$ws.Range("L76").Value = 182000
$ws.Range("N76").Value = -182630

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H79").Value = 186289.25
$ws.Range("J79").Value = 182000
$ws.Range("L79").Value = 182000
$ws.Range("N79").Value = -184184

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 31166.676
$ws.Range("I102").Value = 1402.2
$ws.Range("J102").Value = 113845.78
$ws.Range("K102").Value = 1402.2
$ws.Range("L102").Value = 113845.78
$ws.Range("M102").Value = 219.8
$ws.Range("N102").Value = -117089.78

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I132").Value = 15154807
$ws.Range("J132").Value = 5232
$ws.Range("K132").Value = 45464421
$ws.Range("L132").Value = 15696
$ws.Range("M132").Value = -45461891
$ws.Range("N132").Value = -20756

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2990.389
$ws.Range("I22").Value = 1426
$ws.Range("K22").Value = 1426
$ws.Range("M22").Value = -1131

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2990.389
$ws.Range("I27").Value = 1426
$ws.Range("K27").Value = 1426
$ws.Range("M27").Value = -1319

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5330.1816
$ws.Range("I40").Value = 5029.25
$ws.Range("K40").Value = 5029.25
$ws.Range("M40").Value = -4893.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 12232.538
$ws.Range("J82").Value = 7089.4287
$ws.Range("L82").Value = 7089.4287
$ws.Range("N82").Value = -7811.4287

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 12232.538
$ws.Range("J85").Value = 7089.4287
$ws.Range("L85").Value = 7089.4287
$ws.Range("N85").Value = -9585.4287

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H112").Value = 84248.25
$ws.Range("J112").Value = 84248.25
$ws.Range("L112").Value = 84248.25
$ws.Range("N112").Value = -87202.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 58837270
$ws.Range("I2").Value = 14150.818
$ws.Range("J2").Value = 166679660
$ws.Range("K2").Value = 14150.818
$ws.Range("L2").Value = 166679660
$ws.Range("M2").Value = -14038.818
$ws.Range("N2").Value = -166679884

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5172.2583
$ws.Range("I122").Value = 5319.4644
$ws.Range("K122").Value = 15958.3932
$ws.Range("M122").Value = -13508.3932

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 14985.934
$ws.Range("I136").Value = 27188.75
$ws.Range("K136").Value = 81566.25
$ws.Range("M136").Value = -79016.25
